$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data set gained two new weekly price records.
# One is inserted before the existing row 11, and the other is inserted
# before what becomes row 25 (originally row 24), pushing the rows below
# down accordingly (old 11-23 -> 12-24, old 24-25 -> 26-27).

$ws.Rows.Item(11).Insert()
$ws.Rows.Item(25).Insert()

# New row 11
$ws.Cells.Item(11,1).Value = 8
$ws.Cells.Item(11,2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(11,3).Value = 'Coquimbo'
$d1 = Get-Date -Year 2022 -Month 3 -Day 25 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(11,4).Value = $d1
$ws.Cells.Item(11,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(11,5).Value = 4
$ws.Cells.Item(11,6).Value = 'Fruta'
$ws.Cells.Item(11,7).Value = 100101
$ws.Cells.Item(11,8).Value = 'Berries'
$ws.Cells.Item(11,9).Value = 100101001
$ws.Cells.Item(11,10).Value = 'Arándano (blue)'
$ws.Cells.Item(11,11).Value = 'Sin especificar'
$ws.Cells.Item(11,12).Value = 'Primera'
$ws.Cells.Item(11,13).Value = 100
$ws.Cells.Item(11,14).Value = 6000
$ws.Cells.Item(11,15).Value = 6500
$ws.Cells.Item(11,16).Value = 6250
$ws.Cells.Item(11,17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(11,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(11,19).Value = 3125
$ws.Cells.Item(11,20).Value = 2

# New row 25
$ws.Cells.Item(25,1).Value = 8
$ws.Cells.Item(25,2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(25,3).Value = 'Coquimbo'
$d2 = Get-Date -Year 2022 -Month 3 -Day 24 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(25,4).Value = $d2
$ws.Cells.Item(25,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(25,5).Value = 4
$ws.Cells.Item(25,6).Value = 'Fruta'
$ws.Cells.Item(25,7).Value = 100101
$ws.Cells.Item(25,8).Value = 'Berries'
$ws.Cells.Item(25,9).Value = 100101001
$ws.Cells.Item(25,10).Value = 'Arándano (blue)'
$ws.Cells.Item(25,11).Value = 'Sin especificar'
$ws.Cells.Item(25,12).Value = 'Primera'
$ws.Cells.Item(25,13).Value = 160
$ws.Cells.Item(25,14).Value = 6000
$ws.Cells.Item(25,15).Value = 6500
$ws.Cells.Item(25,16).Value = 6250
$ws.Cells.Item(25,17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(25,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(25,19).Value = 3125
$ws.Cells.Item(25,20).Value = 2

Write-Output "Inserted two new weekly rows; sheet now spans $($ws.UsedRange.Rows.Count) rows."
